# Apply the edits described by the diff: add new G/H sample data rows
# and update the sheet view selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1999
$ws.Range("H2").Formula = "=G2+40001"

$ws.Range("G3").Value = 3999
$ws.Range("H3").Formula = "=G3+40001"

$ws.Range("G4").Value = 2010
$ws.Range("H4").Formula = "=G4+40001"

$ws.Range("G5").Value = 2015
$ws.Range("H5").Formula = "=G5+40001"

# Scroll back to top-left and select the whole second row, matching the
# updated sheetView/selection in the workbook.
$ws.Range("A2:XFD2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
